$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need a temporary Text
# number format so Excel stores them as strings (matching the original
# inlineStr cell type), then the format is reset to Normal so the cell
# keeps its original (default) style index.
$textCells = @("D5","D6","D7","D9","D10","D11","D12","D14","D15","D20","D22","D23","D24","D25","D26","D27","D28","D29","D31","D32","D34","D35","D37","D38","D40","D41","D45","D47","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value2 = "43.056.26"
$ws.Range("E2").Value2 = "  +2.47%  "
$ws.Range("D3").Value2 = "2.303.06"
$ws.Range("E3").Value2 = "  +1.98%  "
$ws.Range("E4").Value2 = "  -0.03%  "
$ws.Range("D5").Value2 = "302.44"
$ws.Range("E5").Value2 = "  +1.23%  "
$ws.Range("D6").Value2 = "99.18"
$ws.Range("E6").Value2 = "  +5.77%  "
$ws.Range("D7").Value2 = "0.507"
$ws.Range("E7").Value2 = "  +1.90%  "
$ws.Range("E8").Value2 = "  -0.04%  "
$ws.Range("D9").Value2 = "0.507"
$ws.Range("E9").Value2 = "  +3.03%  "
$ws.Range("D10").Value2 = "34.38"
$ws.Range("E10").Value2 = "  +4.23%  "
$ws.Range("D11").Value2 = "0.0799"
$ws.Range("E11").Value2 = "  +1.33%  "
$ws.Range("D12").Value2 = "49.19"
$ws.Range("E12").Value2 = "  +3.42%  "
$ws.Range("E13").Value2 = "  +4.18%  "
$ws.Range("D14").Value2 = "18.12"
$ws.Range("E14").Value2 = "  +18.02%  "
$ws.Range("D15").Value2 = "6.81"
$ws.Range("E15").Value2 = "  +2.13%  "
$ws.Range("D16").Value2 = "2.660.50"
$ws.Range("E16").Value2 = "  +1.93%  "
$ws.Range("D17").Value2 = "2.286.73"
$ws.Range("E17").Value2 = "  +1.75%  "
$ws.Range("E18").Value2 = "  +4.59%  "
$ws.Range("D19").Value2 = "42.960.58"
$ws.Range("E19").Value2 = "  +2.22%  "
$ws.Range("D20").Value2 = "12.41"
$ws.Range("E20").Value2 = "  +9.21%  "
$ws.Range("D21").Value2 = "0.0₃0906"
$ws.Range("E21").Value2 = "  +1.72%  "
$ws.Range("D22").Value2 = "6.11"
$ws.Range("E22").Value2 = "  +1.70%  "
$ws.Range("D23").Value2 = "67.96"
$ws.Range("E23").Value2 = "  +2.17%  "
$ws.Range("D24").Value2 = "236.62"
$ws.Range("E24").Value2 = "  +1.35%  "
$ws.Range("D25").Value2 = "2.21"
$ws.Range("E25").Value2 = "  +15.78%  "
$ws.Range("B26").Value2 = "PancakeSwap"
$ws.Range("C26").Value2 = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").Value2 = "2.47"
$ws.Range("E26").Value2 = "  +0.68%  "
$ws.Range("B27").Value2 = "Dai"
$ws.Range("C27").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value2 = "1.00"
$ws.Range("E27").Value2 = "  +0.03%  "
$ws.Range("D28").Value2 = "24.75"
$ws.Range("E28").Value2 = "  +4.50%  "
$ws.Range("D29").Value2 = "168.40"
$ws.Range("E29").Value2 = "  +0.42%  "
$ws.Range("E30").Value2 = "  -9.00%  "
$ws.Range("D31").Value2 = "33.82"
$ws.Range("E31").Value2 = "  +0.79%  "
$ws.Range("D32").Value2 = "9.16"
$ws.Range("E32").Value2 = "  +1.23%  "
$ws.Range("E33").Value2 = "  +0.06%  "
$ws.Range("D34").Value2 = "5.03"
$ws.Range("E34").Value2 = "  +2.02%  "
$ws.Range("D35").Value2 = "4.55"
$ws.Range("E35").Value2 = "  +2.35%  "
$ws.Range("E36").Value2 = "  +3.86%  "
$ws.Range("D37").Value2 = "16.96"
$ws.Range("E37").Value2 = "  +6.25%  "
$ws.Range("D38").Value2 = "0.0701"
$ws.Range("E38").Value2 = "  +1.13%  "
$ws.Range("E39").Value2 = "  +3.76%  "
$ws.Range("B40").Value2 = "ARBITRUM"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value2 = "1.79"
$ws.Range("E40").Value2 = "  +4.71%  "
$ws.Range("B41").Value2 = "LidoDAOToken"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").Value2 = "2.81"
$ws.Range("E41").Value2 = "  +0.77%  "
$ws.Range("E42").Value2 = "  +0.06%  "
$ws.Range("E43").Value2 = "  -2.26%  "
$ws.Range("D44").Value2 = "1.999.60"
$ws.Range("E44").Value2 = "  +2.58%  "
$ws.Range("D45").Value2 = "0.0286"
$ws.Range("E45").Value2 = "  +2.66%  "
$ws.Range("E46").Value2 = "  +4.98%  "
$ws.Range("D47").Value2 = "17.62"
$ws.Range("E47").Value2 = "  +1.57%  "
$ws.Range("D48").Value2 = "2.87"
$ws.Range("E48").Value2 = "  +2.85%  "
$ws.Range("D49").Value2 = "55.54"
$ws.Range("E49").Value2 = "  +6.64%  "
$ws.Range("D50").Value2 = "2.529.26"
$ws.Range("E50").Value2 = "  +1.82%  "
$ws.Range("E51").Value2 = "  +3.11%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
